$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders (column C) and Average (column D) values for the week
$ws.Range("C2").Value = 267
$ws.Range("D2").Value = 231.06

$ws.Range("C3").Value = 243
$ws.Range("D3").Value = 212.6

$ws.Range("C4").Value = 206
$ws.Range("D4").Value = 211.75

$ws.Range("C5").Value = 415
$ws.Range("D5").Value = 239.89

$ws.Range("C6").Value = 213
$ws.Range("D6").Value = 237.5

$ws.Range("C7").Value = 114
$ws.Range("D7").Value = 116.9

$ws.Range("C8").Value = 74
$ws.Range("D8").Value = 96.63

$wb.Save()
